$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 451..467 down to 452..468 (process bottom-up to avoid clobbering
# data before it is read), preserving raw values (Value2 keeps date serials
# as numbers rather than formatted strings) and number formatting (column D
# carries a date number format).
for ($r = 467; $r -ge 451; $r--) {
    for ($c = 1; $c -le 18; $c++) {
        $src = $ws.Cells.Item($r, $c)
        $dst = $ws.Cells.Item($r + 1, $c)
        $dst.Value = $src.Value2()
    }
    # Column D (date) is the only column carrying a non-default style in
    # this sheet; re-apply it explicitly since row 468 is brand new and
    # would otherwise fall back to the default (unstyled) format.
    $ws.Cells.Item($r + 1, 4).NumberFormat = $ws.Cells.Item($r, 4).NumberFormat()
}

# Write the new record into row 451 (a new weekly price observation).
$ws.Cells.Item(451, 1).Value = 4
$ws.Cells.Item(451, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(451, 3).Value = "Los Lagos"
$ws.Cells.Item(451, 4).Value = 45075
$ws.Cells.Item(451, 5).Value = 10
$ws.Cells.Item(451, 6).Value = 100112017
$ws.Cells.Item(451, 7).Value = "Apio"
$ws.Cells.Item(451, 8).Value = "Americana (o)"
$ws.Cells.Item(451, 9).Value = "Primera"
$ws.Cells.Item(451, 10).Value = 25
$ws.Cells.Item(451, 11).Value = 11000
$ws.Cells.Item(451, 12).Value = 11000
$ws.Cells.Item(451, 13).Value = 11000
$ws.Cells.Item(451, 14).Value = "`$/docena de matas"
$ws.Cells.Item(451, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(451, 16).Value = 1833
$ws.Cells.Item(451, 17).Value = 6
$ws.Cells.Item(451, 18).Value = "Hortaliza"
